$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: receive / GAIN / 3 / 1
$ws.Range("A5").Value = "receive"
$ws.Range("E5").Value = "GAIN"
$ws.Range("F5").Value = 3.0
$ws.Range("G5").Value = 1.0

# Row 6: lose / LOSE / 3 / -1
$ws.Range("A6").Value = "lose"
$ws.Range("E6").Value = "LOSE"
$ws.Range("F6").Value = 3.0
$ws.Range("G6").Value = -1.0

# Apply same formatting (style index 3) as the other data rows (rows 2-4),
# matching the sparse column layout (A, E:G only -- no B/C/D cells)
# so we reuse the existing style instead of creating a brand new one.
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

$ws.Range("E4:G4").Copy()
$ws.Range("E5:G6").PasteSpecial(-4122)
